{"js": "// Apply cell-value updates to the single 5-column practice table.\n// Each entry targets an exact (row, col) cell; `old` is asserted defensively\n// so the right cell is touched even though several cells share text.\nconst updates = [\n  { row: 0, col: 0, oldText: \"80\u00f73=26, 2\", newText: \"38\u00f72=19, 0\" },\n  { row: 0, col: 1, oldText: \"76\u00f72=38, 0\", newText: \"51\u00f73=17, 0\" },\n  { row: 0, col: 2, oldText: \"57\u00f74=14, 1\", newText: \"50\u00f73=16, 2\" },\n  { row: 0, col: 3, oldText: \"77\u00f76=12, 5\", newText: \"49\u00f75=9, 4\" },\n  { row: 0, col: 4, oldText: \"67\u00f77=9, 4\", newText: \"59\u00f76=9, 5\" },\n  { row: 4, col: 0, oldText: \"67\u00f77=9, 4\", newText: \"72\u00f73=24, 0\" },\n  { row: 4, col: 1, oldText: \"59\u00f79=6, 5\", newText: \"16\u00f76=2, 4\" },\n  { row: 4, col: 2, oldText: \"82\u00f78=10, 2\", newText: \"12\u00f79=1, 3\" },\n  { row: 4, col: 3, oldText: \"31\u00f75=6, 1\", newText: \"82\u00f72=41, 0\" },\n  { row: 4, col: 4, oldText: \"19\u00f73=6, 1\", newText: \"67\u00f75=13, 2\" },\n  { row: 8, col: 0, oldText: \"64\u00f79=7, 1\", newText: \"45\u00f77=6, 3\" },\n  { row: 8, col: 1, oldText: \"92\u00f74=23, 0\", newText: \"31\u00f79=3, 4\" },\n  { row: 8, col: 2, oldText: \"32\u00f74=8, 0\", newText: \"31\u00f74=7, 3\" },\n  { row: 8, col: 3, oldText: \"96\u00f73=32, 0\", newText: \"97\u00f75=19, 2\" },\n  { row: 8, col: 4, oldText: \"81\u00f77=11, 4\", newText: \"28\u00f78=3, 4\" },\n  { row: 12, col: 0, oldText: \"78\u00f76=13, 0\", newText: \"43\u00f75=8, 3\" },\n  { row: 12, col: 1, oldText: \"68\u00f76=11, 2\", newText: \"81\u00f72=40, 1\" },\n  { row: 12, col: 2, oldText: \"37\u00f72=18, 1\", newText: \"55\u00f76=9, 1\" },\n  { row: 12, col: 3, oldText: \"56\u00f74=14, 0\", newText: \"41\u00f76=6, 5\" },\n  { row: 12, col: 4, oldText: \"94\u00f74=23, 2\", newText: \"53\u00f73=17, 2\" },\n  { row: 16, col: 0, oldText: \"15\u00f77=2, 1\", newText: \"98\u00f77=14, 0\" },\n  { row: 16, col: 1, oldText: \"81\u00f77=11, 4\", newText: \"95\u00f72=47, 1\" },\n  { row: 16, col: 2, oldText: \"77\u00f73=25, 2\", newText: \"28\u00f75=5, 3\" },\n  { row: 16, col: 3, oldText: \"78\u00f77=11, 1\", newText: \"91\u00f72=45, 1\" },\n  { row: 16, col: 4, oldText: \"16\u00f78=2, 0\", newText: \"21\u00f72=10, 1\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nfor (const u of updates) {\n  const current = table.values[u.row][u.col];\n  if (current !== u.oldText) {\n    throw new Error(\n      `Unexpected text at row ${u.row}, col ${u.col}: ` +\n      `expected \"${u.oldText}\" but found \"${current}\"`\n    );\n  }\n  table.getCell(u.row, u.col).value = u.newText;\n}\n\nawait context.sync();", "ps1": "# Apply cell-value updates to the single 5-column practice table.\n# Each entry targets an exact (row, col) cell (1-based, Word COM style);\n# `Old` is asserted defensively so the right cell is touched even though\n# several cells across the table share identical text.\n$updates = @(\n  [PSCustomObject]@{ Row = 1; Col = 1; OldText = \"80\u00f73=26, 2\"; NewText = \"38\u00f72=19, 0\" }\n  [PSCustomObject]@{ Row = 1; Col = 2; OldText = \"76\u00f72=38, 0\"; NewText = \"51\u00f73=17, 0\" }\n  [PSCustomObject]@{ Row = 1; Col = 3; OldText = \"57\u00f74=14, 1\"; NewText = \"50\u00f73=16, 2\" }\n  [PSCustomObject]@{ Row = 1; Col = 4; OldText = \"77\u00f76=12, 5\"; NewText = \"49\u00f75=9, 4\" }\n  [PSCustomObject]@{ Row = 1; Col = 5; OldText = \"67\u00f77=9, 4\"; NewText = \"59\u00f76=9, 5\" }\n  [PSCustomObject]@{ Row = 5; Col = 1; OldText = \"67\u00f77=9, 4\"; NewText = \"72\u00f73=24, 0\" }\n  [PSCustomObject]@{ Row = 5; Col = 2; OldText = \"59\u00f79=6, 5\"; NewText = \"16\u00f76=2, 4\" }\n  [PSCustomObject]@{ Row = 5; Col = 3; OldText = \"82\u00f78=10, 2\"; NewText = \"12\u00f79=1, 3\" }\n  [PSCustomObject]@{ Row = 5; Col = 4; OldText = \"31\u00f75=6, 1\"; NewText = \"82\u00f72=41, 0\" }\n  [PSCustomObject]@{ Row = 5; Col = 5; OldText = \"19\u00f73=6, 1\"; NewText = \"67\u00f75=13, 2\" }\n  [PSCustomObject]@{ Row = 9; Col = 1; OldText = \"64\u00f79=7, 1\"; NewText = \"45\u00f77=6, 3\" }\n  [PSCustomObject]@{ Row = 9; Col = 2; OldText = \"92\u00f74=23, 0\"; NewText = \"31\u00f79=3, 4\" }\n  [PSCustomObject]@{ Row = 9; Col = 3; OldText = \"32\u00f74=8, 0\"; NewText = \"31\u00f74=7, 3\" }\n  [PSCustomObject]@{ Row = 9; Col = 4; OldText = \"96\u00f73=32, 0\"; NewText = \"97\u00f75=19, 2\" }\n  [PSCustomObject]@{ Row = 9; Col = 5; OldText = \"81\u00f77=11, 4\"; NewText = \"28\u00f78=3, 4\" }\n  [PSCustomObject]@{ Row = 13; Col = 1; OldText = \"78\u00f76=13, 0\"; NewText = \"43\u00f75=8, 3\" }\n  [PSCustomObject]@{ Row = 13; Col = 2; OldText = \"68\u00f76=11, 2\"; NewText = \"81\u00f72=40, 1\" }\n  [PSCustomObject]@{ Row = 13; Col = 3; OldText = \"37\u00f72=18, 1\"; NewText = \"55\u00f76=9, 1\" }\n  [PSCustomObject]@{ Row = 13; Col = 4; OldText = \"56\u00f74=14, 0\"; NewText = \"41\u00f76=6, 5\" }\n  [PSCustomObject]@{ Row = 13; Col = 5; OldText = \"94\u00f74=23, 2\"; NewText = \"53\u00f73=17, 2\" }\n  [PSCustomObject]@{ Row = 17; Col = 1; OldText = \"15\u00f77=2, 1\"; NewText = \"98\u00f77=14, 0\" }\n  [PSCustomObject]@{ Row = 17; Col = 2; OldText = \"81\u00f77=11, 4\"; NewText = \"95\u00f72=47, 1\" }\n  [PSCustomObject]@{ Row = 17; Col = 3; OldText = \"77\u00f73=25, 2\"; NewText = \"28\u00f75=5, 3\" }\n  [PSCustomObject]@{ Row = 17; Col = 4; OldText = \"78\u00f77=11, 1\"; NewText = \"91\u00f72=45, 1\" }\n  [PSCustomObject]@{ Row = 17; Col = 5; OldText = \"16\u00f78=2, 0\"; NewText = \"21\u00f72=10, 1\" }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($u in $updates) {\n  $cell = $t.Cell($u.Row, $u.Col)\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $u.OldText) {\n    throw \"Unexpected text at row $($u.Row), col $($u.Col): expected `\"$($u.OldText)`\" but found `\"$current`\"\"\n  }\n  $cell.Range.Text = $u.NewText\n}\n"}
